$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit permutes which species-observation record occupies which row
# position among rows 7-10,12-21 (row 11 is untouched). Two disjoint cycles:
#   7 -> 18 -> 12 -> 9 -> 20 -> 14 -> 13 -> 21 -> 7
#   8 -> 16 -> 19 -> 17 -> 15 -> 10 -> 8
# ("a -> b" meaning row a now holds what row b used to hold).
#
# Strategy: snapshot each of the 14 rows into a scratch row (orig_row + 1000),
# then copy scratch rows back into the permuted destinations, then wipe the
# scratch rows. Work is done in a handful of disjoint column bands (skipping
# the always-empty J/O/X/AF/AH..AS/AU/AV columns) so we neither fabricate
# placeholder cells in columns that must stay untouched, nor leave stale
# values behind in destination cells whose source was blank (Copy() only
# overwrites cells that have content in the source range).
$bands = @(
    @("A", "I"),
    @("K", "N"),
    @("P", "W"),
    @("Y", "AE"),
    @("AG", "AG"),
    @("AT", "AT"),
    @("AW", "AY"),
)

# Step 1: snapshot originals into scratch rows
foreach ($band in $bands) {
    $ws.Range("$($band[0])1007:$($band[1])1007").ClearContents()
    $ws.Range("$($band[0])7:$($band[1])7").Copy($ws.Range("$($band[0])1007:$($band[1])1007"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1008:$($band[1])1008").ClearContents()
    $ws.Range("$($band[0])8:$($band[1])8").Copy($ws.Range("$($band[0])1008:$($band[1])1008"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1009:$($band[1])1009").ClearContents()
    $ws.Range("$($band[0])9:$($band[1])9").Copy($ws.Range("$($band[0])1009:$($band[1])1009"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1010:$($band[1])1010").ClearContents()
    $ws.Range("$($band[0])10:$($band[1])10").Copy($ws.Range("$($band[0])1010:$($band[1])1010"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1012:$($band[1])1012").ClearContents()
    $ws.Range("$($band[0])12:$($band[1])12").Copy($ws.Range("$($band[0])1012:$($band[1])1012"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1013:$($band[1])1013").ClearContents()
    $ws.Range("$($band[0])13:$($band[1])13").Copy($ws.Range("$($band[0])1013:$($band[1])1013"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1014:$($band[1])1014").ClearContents()
    $ws.Range("$($band[0])14:$($band[1])14").Copy($ws.Range("$($band[0])1014:$($band[1])1014"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1015:$($band[1])1015").ClearContents()
    $ws.Range("$($band[0])15:$($band[1])15").Copy($ws.Range("$($band[0])1015:$($band[1])1015"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1016:$($band[1])1016").ClearContents()
    $ws.Range("$($band[0])16:$($band[1])16").Copy($ws.Range("$($band[0])1016:$($band[1])1016"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1017:$($band[1])1017").ClearContents()
    $ws.Range("$($band[0])17:$($band[1])17").Copy($ws.Range("$($band[0])1017:$($band[1])1017"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1018:$($band[1])1018").ClearContents()
    $ws.Range("$($band[0])18:$($band[1])18").Copy($ws.Range("$($band[0])1018:$($band[1])1018"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1019:$($band[1])1019").ClearContents()
    $ws.Range("$($band[0])19:$($band[1])19").Copy($ws.Range("$($band[0])1019:$($band[1])1019"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1020:$($band[1])1020").ClearContents()
    $ws.Range("$($band[0])20:$($band[1])20").Copy($ws.Range("$($band[0])1020:$($band[1])1020"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1021:$($band[1])1021").ClearContents()
    $ws.Range("$($band[0])21:$($band[1])21").Copy($ws.Range("$($band[0])1021:$($band[1])1021"))
}

# Step 2: write scratch rows into the permuted destination rows
foreach ($band in $bands) {
    $ws.Range("$($band[0])7:$($band[1])7").ClearContents()
    $ws.Range("$($band[0])1018:$($band[1])1018").Copy($ws.Range("$($band[0])7:$($band[1])7"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])8:$($band[1])8").ClearContents()
    $ws.Range("$($band[0])1016:$($band[1])1016").Copy($ws.Range("$($band[0])8:$($band[1])8"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])9:$($band[1])9").ClearContents()
    $ws.Range("$($band[0])1020:$($band[1])1020").Copy($ws.Range("$($band[0])9:$($band[1])9"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])10:$($band[1])10").ClearContents()
    $ws.Range("$($band[0])1008:$($band[1])1008").Copy($ws.Range("$($band[0])10:$($band[1])10"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])12:$($band[1])12").ClearContents()
    $ws.Range("$($band[0])1009:$($band[1])1009").Copy($ws.Range("$($band[0])12:$($band[1])12"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])13:$($band[1])13").ClearContents()
    $ws.Range("$($band[0])1021:$($band[1])1021").Copy($ws.Range("$($band[0])13:$($band[1])13"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])14:$($band[1])14").ClearContents()
    $ws.Range("$($band[0])1013:$($band[1])1013").Copy($ws.Range("$($band[0])14:$($band[1])14"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])15:$($band[1])15").ClearContents()
    $ws.Range("$($band[0])1010:$($band[1])1010").Copy($ws.Range("$($band[0])15:$($band[1])15"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])16:$($band[1])16").ClearContents()
    $ws.Range("$($band[0])1019:$($band[1])1019").Copy($ws.Range("$($band[0])16:$($band[1])16"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])17:$($band[1])17").ClearContents()
    $ws.Range("$($band[0])1015:$($band[1])1015").Copy($ws.Range("$($band[0])17:$($band[1])17"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])18:$($band[1])18").ClearContents()
    $ws.Range("$($band[0])1012:$($band[1])1012").Copy($ws.Range("$($band[0])18:$($band[1])18"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])19:$($band[1])19").ClearContents()
    $ws.Range("$($band[0])1017:$($band[1])1017").Copy($ws.Range("$($band[0])19:$($band[1])19"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])20:$($band[1])20").ClearContents()
    $ws.Range("$($band[0])1014:$($band[1])1014").Copy($ws.Range("$($band[0])20:$($band[1])20"))
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])21:$($band[1])21").ClearContents()
    $ws.Range("$($band[0])1007:$($band[1])1007").Copy($ws.Range("$($band[0])21:$($band[1])21"))
}

# Step 3: wipe scratch rows
foreach ($band in $bands) {
    $ws.Range("$($band[0])1007:$($band[1])1007").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1008:$($band[1])1008").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1009:$($band[1])1009").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1010:$($band[1])1010").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1012:$($band[1])1012").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1013:$($band[1])1013").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1014:$($band[1])1014").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1015:$($band[1])1015").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1016:$($band[1])1016").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1017:$($band[1])1017").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1018:$($band[1])1018").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1019:$($band[1])1019").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1020:$($band[1])1020").ClearContents()
}
foreach ($band in $bands) {
    $ws.Range("$($band[0])1021:$($band[1])1021").ClearContents()
}
